# Update the quarterly rial_cumulative income-statement workbook:
#  - drop the two oldest financial-period columns (6m->1399/04, 9m->1399/07)
#  - shift every remaining period one step to the left
#  - append the two newest financial periods (12m->1401/10, 3m->1402/01)
#    together with their published figures
#  - refresh the "publish date" row for the shifted/added columns
#  - minor cosmetic row-height tweaks that came along with the resave

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 8 (B8 = "دوره مالی"): financial-period column headers, D8:M8
# ---------------------------------------------------------------------
$row8 = @(
    "12 ماهه منتهی به 1399/10",
    "3 ماهه منتهی به 1400/01",
    "6 ماهه منتهی به 1400/04",
    "9 ماهه منتهی به 1400/07",
    "12 ماهه منتهی به 1400/10",
    "3 ماهه منتهی به 1401/01",
    "6 ماهه منتهی به 1401/04",
    "9 ماهه منتهی به 1401/07",
    "12 ماهه منتهی به 1401/10",
    "3 ماهه منتهی به 1402/01"
)

# ---------------------------------------------------------------------
# Row 9 (B9 = "تاریخ انتشار"): publish-date column headers, D9:M9
# ---------------------------------------------------------------------
$row9 = @(
    "1401-02-19 (9)",
    "1401-03-01 (2)",
    "1401-07-12 (4)",
    "1401-08-29 (2)",
    "1402-02-13 (9)",
    "1402-02-29 (2)",
    "1401-07-12 (2)",
    "1401-08-29",
    "1402-02-29 (3)",
    "1402-02-29"
)

$col = 4
foreach ($v in $row8) {
    $ws.Cells.Item(8, $col).Value = $v
    $col = $col + 1
}

$col = 4
foreach ($v in $row9) {
    $ws.Cells.Item(9, $col).Value = $v
    $col = $col + 1
}

# ---------------------------------------------------------------------
# Data rows 11-27, columns D:M -- each row shifts two periods to the
# left and gains two new trailing values.
# ---------------------------------------------------------------------
$dataRows = @{
    11 = @(4169376, 1371183, 2653741, 4339949, 6258519, 1875304, 5485485, 8339669, 10301791, 3019253)
    12 = @(-1852457, -752487, -1136959, -1993425, -3031353, -876468, -2131730, -3369967, -4828156, -1284053)
    13 = @(2316919, 618696, 1516782, 2346524, 3227166, 998836, 3353755, 4969702, 5473635, 1735200)
    14 = @(-570671, -203231, -362081, -525482, -768193, -302853, -806640, -1214625, -1370998, -225800)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @(-3654, -4389, -43956, -22003, -100437, -12279, -35775, -17874, -13621, 19526)
    17 = @(1742594, 411076, 1110745, 1799039, 2358536, 683704, 2511340, 3737203, 4089016, 1528926)
    18 = @(-118366, -10159, -28127, -52855, -66371, -4900, -13806, -13352, -15238, 0)
    19 = @(259150, 1410, 14138, 31607, 834047, 46467, 148921, 314917, 1072400, 102205)
    20 = @(1883378, 402327, 1096756, 1777791, 3126212, 725271, 2646455, 4038768, 5146178, 1631131)
    21 = @(-194579, -79326, -163078, -276375, -407103, -108195, -267961, -541819, -425896, -164216)
    22 = @(1688799, 323001, 933678, 1501416, 2719109, 617076, 2378494, 3496949, 4720282, 1466915)
    23 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    24 = @(1688799, 323001, 933678, 1501416, 2719109, 617076, 2378494, 3496949, 4720282, 1466915)
    25 = @(2317, 443, 1281, 2060, 3731, 847, 3264, 4798, 6477, 2013)
    26 = @(728789, 728789, 728789, 728789, 728789, 728789, 728789, 728789, 728789, 728789)
    27 = @(2317, 443, 1281, 2060, 3731, 847, 3264, 4798, 6477, 2013)
}

foreach ($r in $dataRows.Keys) {
    $vals = $dataRows[$r]
    $col = 4
    foreach ($v in $vals) {
        $ws.Cells.Item([int]$r, $col).Value = $v
        $col = $col + 1
    }
}

# ---------------------------------------------------------------------
# Column widths: the two dropped/added periods shift the "wide" (29)
# columns over by one, so re-apply the 28/29 banding across D:M.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 28.1666666666667
$ws.Columns.Item(5).ColumnWidth = 27.1666666666667
$ws.Columns.Item(6).ColumnWidth = 27.1666666666667
$ws.Columns.Item(7).ColumnWidth = 27.1666666666667
$ws.Columns.Item(8).ColumnWidth = 28.1666666666667
$ws.Columns.Item(9).ColumnWidth = 27.1666666666667
$ws.Columns.Item(10).ColumnWidth = 27.1666666666667
$ws.Columns.Item(11).ColumnWidth = 27.1666666666667
$ws.Columns.Item(12).ColumnWidth = 28.1666666666667
$ws.Columns.Item(13).ColumnWidth = 27.1666666666667

# ---------------------------------------------------------------------
# Row-height tweaks carried over by the resave.
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 15.6
$ws.Rows.Item(5).RowHeight = 40.8
$ws.Rows.Item(6).RowHeight = 40.8
$ws.Rows.Item(8).RowHeight = 15.6
